$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at position 451 (shifts existing 451-485 down to 456-490)
$ws.Rows("451:455").Insert()

# Row 451
$ws.Cells.Item(451,1).Value = 5
$ws.Cells.Item(451,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(451,3).Value = "Maule"
$ws.Cells.Item(451,4).Value = 44931
$ws.Cells.Item(451,5).Value = 7
$ws.Cells.Item(451,6).Value = "Fruta"
$ws.Cells.Item(451,7).Value = 100103
$ws.Cells.Item(451,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(451,9).Value = 100103004
$ws.Cells.Item(451,10).Value = "Durazno"
$ws.Cells.Item(451,11).Value = "Kurakata"
$ws.Cells.Item(451,12).Value = "Especial"
$ws.Cells.Item(451,13).Value = 120
$ws.Cells.Item(451,14).Value = 17000
$ws.Cells.Item(451,15).Value = 17000
$ws.Cells.Item(451,16).Value = 17000
$ws.Cells.Item(451,17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(451,18).Value = "Región de O'Higgins"
$ws.Cells.Item(451,19).Value = 1133
$ws.Cells.Item(451,20).Value = 15

# Row 452
$ws.Cells.Item(452,1).Value = 5
$ws.Cells.Item(452,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(452,3).Value = "Maule"
$ws.Cells.Item(452,4).Value = 44931
$ws.Cells.Item(452,5).Value = 7
$ws.Cells.Item(452,6).Value = "Fruta"
$ws.Cells.Item(452,7).Value = 100103
$ws.Cells.Item(452,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(452,9).Value = 100103004
$ws.Cells.Item(452,10).Value = "Durazno"
$ws.Cells.Item(452,11).Value = "Kurakata"
$ws.Cells.Item(452,12).Value = "Primera"
$ws.Cells.Item(452,13).Value = 100
$ws.Cells.Item(452,14).Value = 15000
$ws.Cells.Item(452,15).Value = 15000
$ws.Cells.Item(452,16).Value = 15000
$ws.Cells.Item(452,17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(452,18).Value = "Región de O'Higgins"
$ws.Cells.Item(452,19).Value = 1000
$ws.Cells.Item(452,20).Value = 15

# Row 453
$ws.Cells.Item(453,1).Value = 5
$ws.Cells.Item(453,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(453,3).Value = "Maule"
$ws.Cells.Item(453,4).Value = 44931
$ws.Cells.Item(453,5).Value = 7
$ws.Cells.Item(453,6).Value = "Fruta"
$ws.Cells.Item(453,7).Value = 100103
$ws.Cells.Item(453,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(453,9).Value = 100103004
$ws.Cells.Item(453,10).Value = "Durazno"
$ws.Cells.Item(453,11).Value = "Rich Lady"
$ws.Cells.Item(453,12).Value = "Especial"
$ws.Cells.Item(453,13).Value = 90
$ws.Cells.Item(453,14).Value = 16000
$ws.Cells.Item(453,15).Value = 16000
$ws.Cells.Item(453,16).Value = 16000
$ws.Cells.Item(453,17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(453,18).Value = "Región de O'Higgins"
$ws.Cells.Item(453,19).Value = 1067
$ws.Cells.Item(453,20).Value = 15

# Row 454
$ws.Cells.Item(454,1).Value = 5
$ws.Cells.Item(454,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(454,3).Value = "Maule"
$ws.Cells.Item(454,4).Value = 44931
$ws.Cells.Item(454,5).Value = 7
$ws.Cells.Item(454,6).Value = "Fruta"
$ws.Cells.Item(454,7).Value = 100103
$ws.Cells.Item(454,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(454,9).Value = 100103004
$ws.Cells.Item(454,10).Value = "Durazno"
$ws.Cells.Item(454,11).Value = "Toscana"
$ws.Cells.Item(454,12).Value = "Especial"
$ws.Cells.Item(454,13).Value = 180
$ws.Cells.Item(454,14).Value = 17000
$ws.Cells.Item(454,15).Value = 17000
$ws.Cells.Item(454,16).Value = 17000
$ws.Cells.Item(454,17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(454,18).Value = "Región de O'Higgins"
$ws.Cells.Item(454,19).Value = 1133
$ws.Cells.Item(454,20).Value = 15

# Row 455
$ws.Cells.Item(455,1).Value = 5
$ws.Cells.Item(455,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(455,3).Value = "Maule"
$ws.Cells.Item(455,4).Value = 44931
$ws.Cells.Item(455,5).Value = 7
$ws.Cells.Item(455,6).Value = "Fruta"
$ws.Cells.Item(455,7).Value = 100103
$ws.Cells.Item(455,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(455,9).Value = 100103004
$ws.Cells.Item(455,10).Value = "Durazno"
$ws.Cells.Item(455,11).Value = "Toscana"
$ws.Cells.Item(455,12).Value = "Primera"
$ws.Cells.Item(455,13).Value = 150
$ws.Cells.Item(455,14).Value = 15000
$ws.Cells.Item(455,15).Value = 15000
$ws.Cells.Item(455,16).Value = 15000
$ws.Cells.Item(455,17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(455,18).Value = "Región de O'Higgins"
$ws.Cells.Item(455,19).Value = 1000
$ws.Cells.Item(455,20).Value = 15
